$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 39169
$ws.Cells.Item(2, 4).Value = 56622841
$ws.Cells.Item(3, 3).Value = 93827
$ws.Cells.Item(3, 4).Value = 137518939
$ws.Cells.Item(4, 3).Value = 31991
$ws.Cells.Item(4, 4).Value = 47369138
$ws.Cells.Item(5, 3).Value = 9012
$ws.Cells.Item(5, 4).Value = 13394746
$ws.Cells.Item(6, 3).Value = 2106
$ws.Cells.Item(6, 4).Value = 3130971
$ws.Cells.Item(12, 3).Value = 42553
$ws.Cells.Item(12, 4).Value = 57678152
$ws.Cells.Item(13, 3).Value = 9998
$ws.Cells.Item(13, 4).Value = 14450961
$ws.Cells.Item(14, 3).Value = 26661
$ws.Cells.Item(14, 4).Value = 39082404
$ws.Cells.Item(15, 3).Value = 8505
$ws.Cells.Item(15, 4).Value = 12622978
$ws.Cells.Item(16, 3).Value = 2222
$ws.Cells.Item(16, 4).Value = 3301539
$ws.Cells.Item(20, 3).Value = 10469
$ws.Cells.Item(20, 4).Value = 13843198
$ws.Cells.Item(21, 3).Value = 13803
$ws.Cells.Item(21, 4).Value = 19916328
$ws.Cells.Item(22, 3).Value = 32442
$ws.Cells.Item(22, 4).Value = 47592718
$ws.Cells.Item(23, 3).Value = 10478
$ws.Cells.Item(23, 4).Value = 15572869
$ws.Cells.Item(24, 3).Value = 2716
$ws.Cells.Item(24, 4).Value = 4038771
$ws.Cells.Item(25, 3).Value = 552
$ws.Cells.Item(25, 4).Value = 822092
$ws.Cells.Item(27, 3).Value = 11983
$ws.Cells.Item(27, 4).Value = 15988915
$ws.Cells.Item(28, 3).Value = 7937
$ws.Cells.Item(28, 4).Value = 11485752
$ws.Cells.Item(29, 3).Value = 23163
$ws.Cells.Item(29, 4).Value = 33998287
$ws.Cells.Item(30, 3).Value = 7988
$ws.Cells.Item(30, 4).Value = 11880892
$ws.Cells.Item(31, 3).Value = 2030
$ws.Cells.Item(31, 4).Value = 3029251
$ws.Cells.Item(34, 3).Value = 8544
$ws.Cells.Item(34, 4).Value = 11282927
$ws.Cells.Item(35, 3).Value = 3386
$ws.Cells.Item(35, 4).Value = 4889638
$ws.Cells.Item(36, 3).Value = 8089
$ws.Cells.Item(36, 4).Value = 11815436
$ws.Cells.Item(37, 3).Value = 3250
$ws.Cells.Item(37, 4).Value = 4817961
$ws.Cells.Item(38, 3).Value = 842
$ws.Cells.Item(38, 4).Value = 1254055
$ws.Cells.Item(41, 3).Value = 2551
$ws.Cells.Item(41, 4).Value = 3447067
$ws.Cells.Item(42, 3).Value = 17874
$ws.Cells.Item(42, 4).Value = 25836666
$ws.Cells.Item(43, 3).Value = 52570
$ws.Cells.Item(43, 4).Value = 77046509
$ws.Cells.Item(44, 3).Value = 19405
$ws.Cells.Item(44, 4).Value = 28814623
$ws.Cells.Item(45, 3).Value = 5784
$ws.Cells.Item(45, 4).Value = 8608590
$ws.Cells.Item(46, 3).Value = 1273
$ws.Cells.Item(46, 4).Value = 1900045
$ws.Cells.Item(50, 3).Value = 17246
$ws.Cells.Item(50, 4).Value = 22902417
$ws.Cells.Item(51, 3).Value = 2171
$ws.Cells.Item(51, 4).Value = 3151376
$ws.Cells.Item(52, 3).Value = 7331
$ws.Cells.Item(52, 4).Value = 10771031
$ws.Cells.Item(53, 3).Value = 2455
$ws.Cells.Item(53, 4).Value = 3666572
$ws.Cells.Item(55, 3).Value = 202
$ws.Cells.Item(55, 4).Value = 299226
$ws.Cells.Item(57, 3).Value = 7430
$ws.Cells.Item(57, 4).Value = 10218337
$ws.Cells.Item(58, 3).Value = 1241
$ws.Cells.Item(58, 4).Value = 2179084
$ws.Cells.Item(59, 3).Value = 3024
$ws.Cells.Item(59, 4).Value = 5289247
$ws.Cells.Item(60, 3).Value = 1193
$ws.Cells.Item(60, 4).Value = 2094053
$ws.Cells.Item(64, 3).Value = 1822
$ws.Cells.Item(64, 4).Value = 2946722
$ws.Cells.Item(65, 3).Value = 15962
$ws.Cells.Item(65, 4).Value = 23053147
$ws.Cells.Item(66, 3).Value = 46073
$ws.Cells.Item(66, 4).Value = 67395157
$ws.Cells.Item(67, 3).Value = 16080
$ws.Cells.Item(67, 4).Value = 23891544
$ws.Cells.Item(68, 3).Value = 4688
$ws.Cells.Item(68, 4).Value = 6981788
$ws.Cells.Item(69, 3).Value = 978
$ws.Cells.Item(69, 4).Value = 1453802
$ws.Cells.Item(70, 3).Value = 83
$ws.Cells.Item(70, 4).Value = 121830
$ws.Cells.Item(73, 3).Value = 15508
$ws.Cells.Item(73, 4).Value = 20418869
$ws.Cells.Item(74, 3).Value = 55670
$ws.Cells.Item(74, 4).Value = 80993514
$ws.Cells.Item(75, 3).Value = 154858
$ws.Cells.Item(75, 4).Value = 228090490
$ws.Cells.Item(76, 3).Value = 66551
$ws.Cells.Item(76, 4).Value = 99151049
$ws.Cells.Item(77, 3).Value = 21358
$ws.Cells.Item(77, 4).Value = 31912922
$ws.Cells.Item(78, 3).Value = 5132
$ws.Cells.Item(78, 4).Value = 7666403
$ws.Cells.Item(79, 3).Value = 296
$ws.Cells.Item(79, 4).Value = 439170
$ws.Cells.Item(85, 3).Value = 54510
$ws.Cells.Item(85, 4).Value = 73999555
$ws.Cells.Item(86, 3).Value = 4812
$ws.Cells.Item(86, 4).Value = 6971864
$ws.Cells.Item(87, 3).Value = 11965
$ws.Cells.Item(87, 4).Value = 17574645
$ws.Cells.Item(88, 3).Value = 3985
$ws.Cells.Item(88, 4).Value = 5937458
$ws.Cells.Item(89, 3).Value = 1382
$ws.Cells.Item(89, 4).Value = 2064111
$ws.Cells.Item(90, 3).Value = 299
$ws.Cells.Item(90, 4).Value = 446012
$ws.Cells.Item(93, 3).Value = 5604
$ws.Cells.Item(93, 4).Value = 7529970
$ws.Cells.Item(94, 3).Value = 1675
$ws.Cells.Item(94, 4).Value = 2414699
$ws.Cells.Item(95, 3).Value = 5400
$ws.Cells.Item(95, 4).Value = 7956250
$ws.Cells.Item(96, 3).Value = 2001
$ws.Cells.Item(96, 4).Value = 2979426
$ws.Cells.Item(97, 3).Value = 713
$ws.Cells.Item(97, 4).Value = 1068460
$ws.Cells.Item(101, 3).Value = 3729
$ws.Cells.Item(101, 4).Value = 4942658
$ws.Cells.Item(102, 3).Value = 760
$ws.Cells.Item(102, 4).Value = 1322525
$ws.Cells.Item(103, 3).Value = 473
$ws.Cells.Item(103, 4).Value = 852727
$ws.Cells.Item(104, 3).Value = 172
$ws.Cells.Item(104, 4).Value = 306180
$ws.Cells.Item(105, 3).Value = 59
$ws.Cells.Item(105, 4).Value = 105000
$ws.Cells.Item(107, 3).Value = 11180
$ws.Cells.Item(107, 4).Value = 16211131
$ws.Cells.Item(108, 3).Value = 29944
$ws.Cells.Item(108, 4).Value = 43972354
$ws.Cells.Item(109, 3).Value = 10034
$ws.Cells.Item(109, 4).Value = 14916726
$ws.Cells.Item(110, 3).Value = 2768
$ws.Cells.Item(110, 4).Value = 4127080
$ws.Cells.Item(114, 3).Value = 10068
$ws.Cells.Item(114, 4).Value = 13284109
$ws.Cells.Item(115, 3).Value = 31488
$ws.Cells.Item(115, 4).Value = 45386787
$ws.Cells.Item(116, 3).Value = 67933
$ws.Cells.Item(116, 4).Value = 99390441
$ws.Cells.Item(117, 3).Value = 21847
$ws.Cells.Item(117, 4).Value = 32460828
$ws.Cells.Item(118, 3).Value = 6213
$ws.Cells.Item(118, 4).Value = 9254799
$ws.Cells.Item(119, 3).Value = 1174
$ws.Cells.Item(119, 4).Value = 1754465
$ws.Cells.Item(124, 3).Value = 26497
$ws.Cells.Item(124, 4).Value = 35356461
$ws.Cells.Item(125, 3).Value = 37361
$ws.Cells.Item(125, 4).Value = 53908036
$ws.Cells.Item(126, 3).Value = 79173
$ws.Cells.Item(126, 4).Value = 115744948
$ws.Cells.Item(127, 3).Value = 24463
$ws.Cells.Item(127, 4).Value = 36307847
$ws.Cells.Item(128, 3).Value = 6566
$ws.Cells.Item(128, 4).Value = 9757623
$ws.Cells.Item(129, 3).Value = 1316
$ws.Cells.Item(129, 4).Value = 1957311
$ws.Cells.Item(133, 3).Value = 32669
$ws.Cells.Item(133, 4).Value = 43359448
$ws.Cells.Item(134, 3).Value = 13745
$ws.Cells.Item(134, 4).Value = 19897615
$ws.Cells.Item(135, 3).Value = 33212
$ws.Cells.Item(135, 4).Value = 48771376
$ws.Cells.Item(136, 3).Value = 11744
$ws.Cells.Item(136, 4).Value = 17449316
$ws.Cells.Item(137, 3).Value = 3056
$ws.Cells.Item(137, 4).Value = 4555241
$ws.Cells.Item(138, 3).Value = 520
$ws.Cells.Item(138, 4).Value = 773990
$ws.Cells.Item(141, 3).Value = 11107
$ws.Cells.Item(141, 4).Value = 14804113
$ws.Cells.Item(142, 3).Value = 36524
$ws.Cells.Item(142, 4).Value = 52752943
$ws.Cells.Item(143, 3).Value = 84150
$ws.Cells.Item(143, 4).Value = 123275408
$ws.Cells.Item(144, 3).Value = 25048
$ws.Cells.Item(144, 4).Value = 37212458
$ws.Cells.Item(145, 3).Value = 6574
$ws.Cells.Item(145, 4).Value = 9809496
$ws.Cells.Item(149, 3).Value = 30085
$ws.Cells.Item(149, 4).Value = 40549891
